# Fruta / hortaliza, semanal
#
# Insert two new price-record rows (row 126 and row 127, "Artic Snow"
# variety) above the existing "Nectar Crest" block in the Nectarín sheet.
# All rows currently at 126..211 shift down two positions to 128..213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting two blank rows at 126:127.
# (Mirrors Excel's normal Rows.Insert() behaviour: formatting - e.g. the
# date NumberFormat style carried on column D - is inherited from the row
# above, same as a manual "Insert Copied/Sheet Rows" in the UI.)
$ws.Rows("126:127").Insert()

# Populate the newly-inserted row 126 ("Artic Snow" / "Especial").
$ws.Range("A126").Value = 2
$ws.Range("B126").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C126").Value = "Coquimbo"
$ws.Range("D126").Value = 44651
$ws.Range("E126").Value = 4
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100103
$ws.Range("H126").Value = "Frutos de hueso (carozo)"
$ws.Range("I126").Value = 100103006
$ws.Range("J126").Value = "Nectarín"
$ws.Range("K126").Value = "Artic Snow"
$ws.Range("L126").Value = "Especial"
$ws.Range("M126").Value = 10
$ws.Range("N126").Value = 430000
$ws.Range("O126").Value = 440000
$ws.Range("P126").Value = 435000
$ws.Range("Q126").Value = "$/bins (420 kilos)"
$ws.Range("R126").Value = "Región de O'Higgins"
$ws.Range("S126").Value = 1036
$ws.Range("T126").Value = 420

# Populate the newly-inserted row 127 ("Artic Snow" / "Primera").
$ws.Range("A127").Value = 2
$ws.Range("B127").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C127").Value = "Coquimbo"
$ws.Range("D127").Value = 44651
$ws.Range("E127").Value = 4
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100103
$ws.Range("H127").Value = "Frutos de hueso (carozo)"
$ws.Range("I127").Value = 100103006
$ws.Range("J127").Value = "Nectarín"
$ws.Range("K127").Value = "Artic Snow"
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 16
$ws.Range("N127").Value = 380000
$ws.Range("O127").Value = 390000
$ws.Range("P127").Value = 385000
$ws.Range("Q127").Value = "$/bins (420 kilos)"
$ws.Range("R127").Value = "Región de O'Higgins"
$ws.Range("S127").Value = 917
$ws.Range("T127").Value = 420
